$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-70 down to 31-71.
$ws.Rows(30).Insert()

# Populate the newly inserted row 30 with the new weekly data point.
$ws.Cells.Item(30, 1).Value = 11
$ws.Cells.Item(30, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(30, 3).Value = "Bíobío"
$ws.Cells.Item(30, 4).Value = 45219
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(30, 6).Value = 100112026
$ws.Cells.Item(30, 7).Value = "Haba"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 50
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 10000
$ws.Cells.Item(30, 13).Value = 10000
$ws.Cells.Item(30, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(30, 16).Value = 400
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
